$wb = $excel.ActiveWorkbook

$oldGuid = "b0cd265a-9723-43ea-9ba4-e888068ac875"
$newGuid = "d8737fb6-aa9c-452c-afa9-e35f63e21a94"

$oldZhHash = "9b7078a71b0291b1a4908e7547d7532fea9a015f"
$newZhHash = "07e02ea9d1f1514fa21061e889d9b7927b342ac4"

# Overview sheet
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Hyperlinks.Item(1).TextToDisplay = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-09-01 15:22:46"

# zh-cn sheet
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = "$newGuid.md"
$wsZh.Hyperlinks.Item(1).TextToDisplay = "$newGuid.md"
$wsZh.Range("G2").Value = "$newGuid.$newZhHash.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-09-01 15:22:40"

# de-de sheet
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = "$newGuid.md"
$wsDe.Hyperlinks.Item(1).TextToDisplay = "$newGuid.md"
$wsDe.Range("H2").Value = "$newGuid.$newZhHash.de-de.xlf"
